$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the minimum-age / observation & process error input cells.
# All other cells on the sheet (A2, B2, B6, B9, B10) are formulas that
# depend on these inputs and will recalculate automatically.
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 1

$excel.Calculate()
